$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46050
$ws.Range("B2").Value = 3.78
$ws.Range("C2").Value = 3.74
$ws.Range("D2").Value = 2.17
$ws.Range("E2").Value = 1.98
$ws.Range("F2").Value = 2.63
$ws.Range("G2").Value = 3.53
$ws.Range("H2").Value = 11.85
$ws.Range("I2").Value = 32.23
$ws.Range("J2").Value = 57.94
$ws.Range("K2").Value = 78.14
$ws.Range("L2").Value = 46.68
$ws.Range("M2").Value = 22.41
$ws.Range("N2").Value = 12.78
$ws.Range("O2").Value = 8.09
$ws.Range("P2").Value = 3.78
$ws.Range("Q2").Value = 3.78
$ws.Range("R2").Value = 7.46
$ws.Range("S2").Value = 23.22
$ws.Range("T2").Value = 48.86
$ws.Range("U2").Value = 75.95999999999999
$ws.Range("V2").Value = 86.28
$ws.Range("W2").Value = 73.59999999999999
$ws.Range("X2").Value = 34.29
$ws.Range("Y2").Value = 13.91
$ws.Range("Z2").Value = 27.46
$ws.Range("AB2").Value = 52.02
$ws.Range("AD2").Value = 79.94
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 68.04000000000001
